# Update gh-pages output data (展览 / 全部类型 sheets) with refreshed
# "想去人数" (F) and "最低票价" (G) values.

$wb = $excel.ActiveWorkbook

function Update-Sheet {
    param($ws)

    $ws.Range("F5").Value = 15463
    $ws.Range("F8").Value = 693
    $ws.Range("F9").Value = 15370
    $ws.Range("G9").Value = 19.9
    $ws.Range("F11").Value = 8966
    $ws.Range("F12").Value = 366
}

# 展览 sheet (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
Update-Sheet $ws1
$ws1.Range("F32").Value = 56
$ws1.Range("F34").Value = 244
$ws1.Range("G35").Value = "不可售"
$ws1.Range("F36").Value = 444
$ws1.Range("F38").Value = 5489

# 全部类型 sheet (sheet4) - same underlying data, rows offset by +2 from row 32 on
$ws4 = $wb.Worksheets.Item("全部类型")
Update-Sheet $ws4
$ws4.Range("F34").Value = 56
$ws4.Range("F36").Value = 244
$ws4.Range("G37").Value = "不可售"
$ws4.Range("F38").Value = 444
$ws4.Range("F40").Value = 5489
